$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Boundaries")

# Make sure this is the active sheet (it already is tabSelected in the file).
$ws.Activate()

# The "15" constant used to offset column F dropped to "14.9" for the whole
# block (G2 standalone formula, plus the G3:G8 shared formula that every
# other row in the block inherits). Downstream columns I/K/L recompute
# automatically since they reference G.
$ws.Range("G2").Formula = "=-F2+14.9"
$ws.Range("G3:G8").Formula = "=-F3+14.9"

# Selection moved from K3 to G11 (and the frozen/scrolled topLeftCell="B1"
# view setting falls away once the selection is simply re-pointed).
$ws.Range("G11").Select() | Out-Null
